# "Generate Report for Handoff"
# The b.md entry has been handed off again: its status moves from
# "Handed back: in sync with en-US" to "Ready for handoff", a new xlf
# handoff package was generated for each locale, and a warning about a
# stale handback version is recorded in the Error Detail column.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---- Overview sheet: row 3 is the b.md file ----
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-11-03 19:19:04"

# ---- zh-cn sheet: row 3 is the b.md file ----
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-11-03 19:18:50"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/387adb36c2bbddf9d5c16425555f37b1802ae935/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/45e3760ef7ee11db0528e7f8167af883d17f9da8/e2e/b.md."

# ---- de-de sheet: row 3 is the b.md file ----
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-11-03 19:19:04"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/387adb36c2bbddf9d5c16425555f37b1802ae935/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/45e3760ef7ee11db0528e7f8167af883d17f9da8/e2e/b.md."

# The Error Detail column now holds a long message, so widen it to 40
# characters on both locale sheets (matches the column width bump seen
# in the workbook XML).
$zhcn.Columns.Item(16).ColumnWidth = 39.16666666666667
$dede.Columns.Item(16).ColumnWidth = 39.16666666666667
